# Update 想去人数 (interest counts) in column F across all sheets
# per "output generated at 456a3b4" data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 2154
$ws.Range("F6").Value = 567
$ws.Range("F13").Value = 666659
$ws.Range("F14").Value = 1663
$ws.Range("F15").Value = 560
$ws.Range("F16").Value = 1482
$ws.Range("F19").Value = 1296
$ws.Range("F20").Value = 2286
$ws.Range("F22").Value = 2713
$ws.Range("F23").Value = 1572
$ws.Range("F24").Value = 861
$ws.Range("F27").Value = 1097
$ws.Range("F28").Value = 531
$ws.Range("F29").Value = 531
$ws.Range("F36").Value = 1313
$ws.Range("F37").Value = 2787
$ws.Range("F39").Value = 1153
$ws.Range("F44").Value = 1003
$ws.Range("F45").Value = 3170

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 112
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 144822
$ws.Range("F23").Value = 85
$ws.Range("F26").Value = 609
$ws.Range("F31").Value = 366
$ws.Range("F38").Value = 217

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F7").Value = 831
$ws.Range("F10").Value = 1612
$ws.Range("F11").Value = 127
$ws.Range("F12").Value = 1993

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 831
$ws.Range("F6").Value = 1612
$ws.Range("F8").Value = 2154
$ws.Range("F9").Value = 127
$ws.Range("F10").Value = 1993
$ws.Range("F12").Value = 567
$ws.Range("F17").Value = 666660
$ws.Range("F18").Value = 112
$ws.Range("F20").Value = 1663
$ws.Range("F21").Value = 144823
$ws.Range("F22").Value = 1482
$ws.Range("F25").Value = 1296
$ws.Range("F26").Value = 2286
$ws.Range("F28").Value = 2713
$ws.Range("F29").Value = 1572
$ws.Range("F30").Value = 861
$ws.Range("F35").Value = 1097
$ws.Range("F36").Value = 531
$ws.Range("F38").Value = 85
$ws.Range("F40").Value = 1313
$ws.Range("F41").Value = 2788
$ws.Range("F43").Value = 1153
$ws.Range("F44").Value = 366
$ws.Range("F48").Value = 1003
$ws.Range("F49").Value = 3170
$ws.Range("F50").Value = 217
